$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.196.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.737.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.735.59'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.54'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.479'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000252'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.356.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.735.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.258.39'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('E18').Value = '  -2.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '498.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.718'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000133'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.112'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.345'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.137'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.22%  '
$ws.Range('E40').Value = '  +12.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '441.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.05'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.944.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0358'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '138.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.05%  '
